$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "36.535.99"
$ws.Range("E2").Value = "  +0.03%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.941.36"
$ws.Range("E3").Value = "  -3.34%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.43%  "

# Row 7 - Solana (was USDC)
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.37"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -8.66%  "

# Row 8 - USDC (was Solana)
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.365"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.00%  "

# Row 10 - OKB
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.80%  "

# Row 11 - Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0835"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.08%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.02%  "

# Row 13 - Polygon
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.823"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.72%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.31%  "

# Row 15 - Wrapped liquid staked Ether 2.0
$ws.Range("D15").Value = "2.224.22"
$ws.Range("E15").Value = "  -3.59%  "

# Row 16 - Chainlink
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.01%  "

# Row 17 - Polkadot
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.67%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "1.941.34"
$ws.Range("E18").Value = "  -3.58%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "36.417.96"
$ws.Range("E19").Value = "  -0.27%  "

# Row 20 - Litecoin
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.32%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0867"
$ws.Range("E21").Value = "  -1.49%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.77%  "

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.07%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.01%  "

# Row 25 - PancakeSwap
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.72%  "

# Row 26 - Toncoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.20%  "

# Row 27 - Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.12%  "

# Row 28 - Monero
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.03%  "

# Row 29 - EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.74%  "

# Row 30 - Kaspa
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.124"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.88%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  -3.09%  "

# Row 32 - ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.10%  "

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.78%  "

# Row 34 - Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0628"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.33%  "

# Row 35 - InternetComputer (DFINITY)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.35%  "

# Row 36 - BinanceUSD
$ws.Range("E36").Value = "  -0.23%  "

# Row 37 - THORChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.98%  "

# Row 38 - WEMIXToken
$ws.Range("E38").Value = "  -3.57%  "

# Row 39 - LidoDAOToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.36%  "

# Row 40 - RenderToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.30%  "

# Row 41 - Cronos
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0972"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.49%  "

# Row 42 - HuobiToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.07%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  -7.01%  "

# Row 44 - VeChain
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0208"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.00%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.37%  "

# Row 46 - ARBITRUM
$ws.Range("E46").Value = "  -8.02%  "

# Row 47 - Maker
$ws.Range("D47").Value = "1.346.86"
$ws.Range("E47").Value = "  -1.40%  "

# Row 48 - Aave
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.34%  "

# Row 49 - FraxShare
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.64%  "

# Row 50 - MXToken
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.46%  "

# Row 51 - MultiversX
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.37%  "
